# services/cards.xlsx fix:
#  - row 6 ("обычный воробей") had placeholder "-" values for its
#    ability type (ТИП) and ability power (сила способности) columns.
#    Update the type to "deathrattle" and the power to the numeric value 0.
#  - the row's height was an oversized 33pt (leftover from the placeholder
#    text); shrink it back down to the standard 14.25pt row height.
#  - move the sheet's active-cell selection to C26 (cursor left where the
#    user had last clicked while editing).
#
# Removing the now-unused "-" shared string and renumbering every other
# shared-string reference is handled automatically by Excel when the
# workbook is saved, so no manual shared-string bookkeeping is needed here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "deathrattle"
$ws.Range("D6").Value = 0

$ws.Rows.Item(6).RowHeight = 14.25

$ws.Range("C26").Select()
